# Insert a new working-hours entry before the "Total" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Total" row currently sits at row 6. Move it down to row 7 (copying
# its cells, formulas and formatting verbatim) so row 6 is free for a new
# data entry, just like rows 2-5. We deliberately avoid Rows.Insert() here
# since it mints a brand-new (unused) style entry in the workbook's style
# table as a side effect of the automatic row-formatting inheritance.
$ws.Range("A6:F6").Copy($ws.Range("A7:F7"))
$excel.CutCopyMode = $false

# New data row (row 6): Date / From / To / Hours / Hour Rate / Bill
$ws.Cells.Item(6, 1).Value = 45268
$ws.Cells.Item(6, 2).Value = 0.708333333333333
$ws.Cells.Item(6, 3).Value = 0.875
$ws.Cells.Item(6, 4).Formula = "=(C6<B6)+C6-B6"
$ws.Cells.Item(6, 5).Value = 10
$ws.Cells.Item(6, 6).Formula = "=(D6*24)*E6"

# Copy the formatting of the previous data row (row 5) onto the new row (6)
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the Total row (now row 7) SUM ranges to include the new row 6
$ws.Cells.Item(7, 4).Formula = "=SUM(D2:D6)"
$ws.Cells.Item(7, 6).Formula = "=SUM(F2:F6)"

$ws.Range("I8").Select()
